$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row: B1 becomes "Comment"
$ws.Range("B1").Value = "Comment"

# Clear out column C content (Second Comment / N/A) - the column width
# definition stays in place, only the cell contents are removed.
$ws.Range("C1:C2").ClearContents()

# Select B2 as the active cell (matches final selection in workbook)
$ws.Range("B2").Select()
